$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "custom data label" header columns (L1:P1) to conform to the
# new naming convention -- SQLite-safe column names are stored as the
# field's localized display label instead of the old English identifier.
$ws.Range("L1").Value = "कस्टम डेटा लेबल के साथ काउंटर"
$ws.Range("M1").Value = "कस्टम डेटा लेबल के साथ काउंटरMarkers"
$ws.Range("P1").Value = "カスタムデータラベルのフラグ"
$ws.Range("O1").Value = "注意使用自定义数据标签"
$ws.Range("N1").Value = 'Bandera con"etiqueta"de"datos personalizada'

# Move the active selection from N3 to N1.
$ws.Range("N1").Select()
